$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("C2").Value = 0.9938
$ws.Range("D2").Value = 0.8242
$ws.Range("E2").Value = 0.0004
$ws.Range("G2").Value = 0.508
$ws.Range("H2").Value = 0.9942

# Row 3 updates
$ws.Range("E3").Value = 0.9996
$ws.Range("F3").Value = 0.9998
